$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.697.22'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '2.508.78'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.49'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.04'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('D9').Value = '2.507.86'
$ws.Range('E9').Value = '  -0.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.162'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.02%  '
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('E12').Value = '  +6.22%  '
$ws.Range('E13').Value = '  +1.60%  '
$ws.Range('D14').Value = '2.980.70'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('D16').Value = '69.469.58'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.88'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '2.517.20'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '351.35'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.92'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.65%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.15'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.79%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.95'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.65%  '
$ws.Range('D28').Value = '2.666.00'
$ws.Range('E28').Value = '  +0.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('D30').Value = '0.0₃0893'
$ws.Range('E30').Value = '  -1.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.89'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '461.72'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.05%  '
$ws.Range('E33').Value = '  -6.44%  '
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.13'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.61%  '
$ws.Range('E37').Value = '  +1.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.09'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.54'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.57%  '
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.319'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.60'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.95%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.22'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.10'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.66'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.522'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.64%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.48'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.95%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0734'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('E51').Value = '  +3.22%  '
